# TC05 Canine workbook update:
#   - add "CypherOutput_Message" sheet (copy of the existing "Message" sheet)
#   - add "StatOutput" sheet with the summary counts
#   - add "StatOutput_Message" sheet (two stacked copies of the "Message"
#     sheet content; the second copy's Cypher query is the stats/count query)

$wb = $excel.ActiveWorkbook

$cypherOutput = $wb.Worksheets.Item("CypherOutput")
$message = $wb.Worksheets.Item("Message")

# The second Cypher query text (the stats / counts query) used on the
# StatOutput_Message sheet - differs from the query already stored in the
# Message sheet (which is reused, unmodified, for CypherOutput_Message).
$statsCypher = @'
MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE demo.sex IN ['Male Phenotype']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study
'@

# ---------------------------------------------------------------------
# 1) CypherOutput_Message - a straight copy of the Message sheet content
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$cypherOutputMessage = $wb.Worksheets.Add($null, $lastSheet)
$cypherOutputMessage.Name = "CypherOutput_Message"

$message.Range("A1:A10").Copy()
$cypherOutputMessage.Range("A1").PasteSpecial(-4163)

# ---------------------------------------------------------------------
# 2) StatOutput - summary counts table
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$statOutput = $wb.Worksheets.Add($null, $lastSheet)
$statOutput.Name = "StatOutput"

$statOutput.Range("A1").Value = "number_of_files"
$statOutput.Range("B1").Value = "number_of_sample"
$statOutput.Range("C1").Value = "number_of_cases"
$statOutput.Range("D1").Value = "number_of_study"

# The counts look like numbers but the source workbook stores them as plain
# text (shared strings), so write them through a text formula and paste the
# computed value back in - this keeps the numeric-looking text as text
# instead of Excel auto-converting it to a real number.
$statOutput.Range("Z1").Formula = '="0"'
$statOutput.Range("Z2").Formula = '="0"'
$statOutput.Range("Z3").Formula = '="8"'
$statOutput.Range("Z4").Formula = '="1"'

$statOutput.Range("Z1").Copy()
$statOutput.Range("A2").PasteSpecial(-4163)
$statOutput.Range("Z2").Copy()
$statOutput.Range("B2").PasteSpecial(-4163)
$statOutput.Range("Z3").Copy()
$statOutput.Range("C2").PasteSpecial(-4163)
$statOutput.Range("Z4").Copy()
$statOutput.Range("D2").PasteSpecial(-4163)

$statOutput.Range("Z1:Z4").ClearContents()

# ---------------------------------------------------------------------
# 3) StatOutput_Message - Message content twice; 2nd block's Cypher query
#    is the stats/count query above.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$statOutputMessage = $wb.Worksheets.Add($null, $lastSheet)
$statOutputMessage.Name = "StatOutput_Message"

$message.Range("A1:A10").Copy()
$statOutputMessage.Range("A1").PasteSpecial(-4163)
$message.Range("A1:A10").Copy()
$statOutputMessage.Range("A11").PasteSpecial(-4163)

# Row 18 (8th row of the 2nd block) holds the Cypher query - replace it with
# the stats query text. (Plain text, not numeric-looking, so a direct
# .Value assignment already lands as a shared string.)
$statOutputMessage.Range("A18").Value = $statsCypher

# Keep the originally-active sheet selected, same as the source workbook.
$cypherOutput.Activate()
